$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the "estado de cuenta" detail rows (16-19) by Periodo Mora, then by
# worker, and update the corresponding "Valor Mora" amounts to match the
# refreshed source database.
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1235038025"
$ws.Range("D16").Value = "RICARDO MARIO JIMENEZ RESTREPO"
$ws.Range("E16").Value = "2409"
$ws.Range("F16").Value = 52000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "13541643"
$ws.Range("D17").Value = "PEDRO ALONSO HERNANDEZ ROMERO"
$ws.Range("E17").Value = "2409"
$ws.Range("F17").Value = 52000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1235038025"
$ws.Range("D18").Value = "RICARDO MARIO JIMENEZ RESTREPO"
$ws.Range("E18").Value = "2410"
$ws.Range("F18").Value = 15600

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "13541643"
$ws.Range("D19").Value = "PEDRO ALONSO HERNANDEZ ROMERO"
$ws.Range("E19").Value = "2410"
$ws.Range("F19").Value = 15600
